$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells stay text (several values look numeric, e.g. "17.60")
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "58.061.66"
$ws.Range("E2").Value = "  +2.49%  "
$ws.Range("D3").Value = "2.350.18"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "545.13"
$ws.Range("E5").Value = "  +5.88%  "
$ws.Range("D6").Value = "134.87"
$ws.Range("E6").Value = "  +2.24%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").Value = "0.537"
$ws.Range("E8").Value = "  +0.66%  "
$ws.Range("D9").Value = "2.346.46"
$ws.Range("E9").Value = "  +0.78%  "
$ws.Range("E10").Value = "  +1.30%  "
$ws.Range("E11").Value = "  +1.15%  "
$ws.Range("D12").Value = "5.42"
$ws.Range("E12").Value = "  +3.43%  "
$ws.Range("D13").Value = "0.358"
$ws.Range("E13").Value = "  +6.35%  "
$ws.Range("D14").Value = "2.767.00"
$ws.Range("E14").Value = "  +1.20%  "
$ws.Range("D15").Value = "23.57"
$ws.Range("E15").Value = "  -0.19%  "
$ws.Range("D16").Value = "58.032.85"
$ws.Range("E16").Value = "  +2.49%  "
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("D18").Value = "2.336.11"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("E19").Value = "  +2.60%  "
$ws.Range("D20").Value = "334.51"
$ws.Range("E20").Value = "  +1.91%  "
$ws.Range("E21").Value = "  +1.45%  "
$ws.Range("D22").Value = "6.71"
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "61.69"
$ws.Range("E24").Value = "  +0.74%  "
$ws.Range("E25").Value = "  +2.62%  "
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("E27").Value = "  -1.89%  "
$ws.Range("E28").Value = "  +7.53%  "
$ws.Range("D29").Value = "1.77"
$ws.Range("E29").Value = "  +5.36%  "
$ws.Range("D30").Value = "170.38"
$ws.Range("E30").Value = "  +1.77%  "
$ws.Range("D31").Value = "0.0₃0731"
$ws.Range("E31").Value = "  +1.50%  "
$ws.Range("D32").Value = "6.13"
$ws.Range("E32").Value = "  +0.53%  "
$ws.Range("D34").Value = "18.45"
$ws.Range("E34").Value = "  +0.87%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  +0.27%  "
$ws.Range("E37").Value = "  +6.27%  "
$ws.Range("E38").Value = "  +1.79%  "
$ws.Range("E39").Value = "  +4.44%  "
$ws.Range("D40").Value = "39.36"
$ws.Range("E40").Value = "  +1.94%  "
$ws.Range("D41").Value = "148.41"
$ws.Range("E41").Value = "  -0.64%  "
$ws.Range("E42").Value = "  +1.30%  "
$ws.Range("D43").Value = "287.35"
$ws.Range("E43").Value = "  +3.95%  "
$ws.Range("E44").Value = "  +1.00%  "
$ws.Range("D45").Value = "19.27"
$ws.Range("E45").Value = "  +5.51%  "
$ws.Range("D46").Value = "0.0927"
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("D47").Value = "0.0504"
$ws.Range("E47").Value = "  +1.95%  "
$ws.Range("D48").Value = "0.563"
$ws.Range("E48").Value = "  +1.38%  "
$ws.Range("D49").Value = "0.0218"
$ws.Range("E49").Value = "  +1.62%  "
$ws.Range("D50").Value = "17.60"
$ws.Range("E50").Value = "  +2.89%  "
$ws.Range("E51").Value = "  +0.54%  "
